# Apply the edit described by the diff:
# 1. For each year group of 4 rows (A/B/C/D pattern starting at row 2),
#    swap the contents (columns A-E) of the "B" row and the "C" row.
# 2. Delete columns F and G (the "十种有色金属产销率" and "十种有色金属销售量" columns),
#    which duplicated data now covered by columns B and E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 81
$firstGroupRow = 2
$groupSize = 4
$lastCol = 5   # columns A..E

for ($groupStart = $firstGroupRow; $groupStart -le $lastDataRow; $groupStart += $groupSize) {
    $rowB = $groupStart + 1
    $rowC = $groupStart + 2

    for ($col = 1; $col -le $lastCol; $col++) {
        $cellB = $ws.Cells.Item($rowB, $col)
        $cellC = $ws.Cells.Item($rowC, $col)

        $valB = $cellB.Value()
        $valC = $cellC.Value()

        $cellB.Value = $valC
        $cellC.Value = $valB
    }
}

# Remove the redundant F and G columns entirely.
$ws.Range("F1:G$lastDataRow").Delete()
